$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I11").Value = 3875.2
$ws.Range("M11").Value = -3735.2
$ws.Range("H11").Value = 3875.2
$ws.Range("K11").Value = 3875.2
$ws.Range("I18").Value = 456.8
$ws.Range("K18").Value = 456.8
$ws.Range("M18").Value = -172.8
$ws.Range("H18").Value = 4081
$ws.Range("J19").Value = 1282.4546
$ws.Range("L19").Value = 1282.4546
$ws.Range("H19").Value = 1255.5238
$ws.Range("N19").Value = -1632.4546
$ws.Range("N43").Value = -1373116.9
$ws.Range("H43").Value = 698062.5
$ws.Range("L43").Value = 1372978.9
$ws.Range("J43").Value = 1372978.9
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("J64").Value = 7981.6665
$ws.Range("L64").Value = 7981.6665
$ws.Range("N64").Value = -8477.666499999999
$ws.Range("H64").Value = 7236.25
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("H67").Value = 7236.25
$ws.Range("J67").Value = 7981.6665
$ws.Range("L67").Value = 7981.6665
$ws.Range("N67").Value = -9697.666499999999
$ws.Range("H86").Value = 183763790
$ws.Range("K86").Value = 166668720
$ws.Range("I86").Value = 166668720
$ws.Range("M86").Value = -166667597
$ws.Range("H89").Value = 183763790
$ws.Range("M89").Value = -833337984
$ws.Range("I89").Value = 166668720
$ws.Range("K89").Value = 833343600
$ws.Range("J98").Value = 1413.8
$ws.Range("I98").Value = 2089.487
$ws.Range("L98").Value = 1413.8
$ws.Range("N98").Value = -4409.8
$ws.Range("K98").Value = 2089.487
$ws.Range("M98").Value = -591.4870000000001
$ws.Range("H98").Value = 1951.5918
$ws.Range("H111").Value = 13891330
$ws.Range("I111").Value = 20836104
$ws.Range("K111").Value = 62508312
$ws.Range("L111").Value = 5350.0002
$ws.Range("N111").Value = -11484.0002
$ws.Range("M111").Value = -62505245
$ws.Range("J111").Value = 1783.3334
$ws.Range("L115").Value = 3000
$ws.Range("J115").Value = 1000
$ws.Range("H115").Value = 637.5
$ws.Range("N115").Value = -6134
$ws.Range("L122").Value = 4241.4
$ws.Range("K122").Value = 6268.461
$ws.Range("H122").Value = 1951.5918
$ws.Range("I122").Value = 2089.487
$ws.Range("N122").Value = -9141.4
$ws.Range("J122").Value = 1413.8
$ws.Range("M122").Value = -3818.461
$ws.Range("N125").Value = -58885.125
$ws.Range("L125").Value = 53965.125
$ws.Range("J125").Value = 5996.125
$ws.Range("H125").Value = 26319454
$ws.Range("J133").Value = 78780
$ws.Range("L133").Value = 78780
$ws.Range("N133").Value = -88900
$ws.Range("H133").Value = 78780
$ws.Range("N137").Value = -11022
$ws.Range("I137").Value = 1468.25
$ws.Range("M137").Value = -1854.75
$ws.Range("H137").Value = 1569.4
$ws.Range("J137").Value = 1974
$ws.Range("K137").Value = 4404.75
$ws.Range("L137").Value = 5922
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 7359736.5
$ws.Range("M32").Value = -7359449.5
$ws.Range("H32").Value = 6258476
$ws.Range("K32").Value = 7359736.5
$ws.Range("H45").Value = 4166
$ws.Range("M45").Value = -1462.2727
$ws.Range("K45").Value = 1839.2727
$ws.Range("I45").Value = 1839.2727
$ws.Range("I61").Value = 3072.5334
$ws.Range("K61").Value = 3072.5334
$ws.Range("M61").Value = -2860.5334
$ws.Range("H61").Value = 7621.121
$ws.Range("I74").Value = 46185
$ws.Range("K74").Value = 46185
$ws.Range("H74").Value = 33178.44
$ws.Range("M74").Value = -45311
$ws.Range("K77").Value = 230925
$ws.Range("I77").Value = 46185
$ws.Range("H77").Value = 33178.44
$ws.Range("M77").Value = -226557
$ws.Range("I110").Value = 414.26666
$ws.Range("H110").Value = 14493312
$ws.Range("M110").Value = 1630.73334
$ws.Range("K110").Value = 414.26666
$ws.Range("I132").Value = 1856393.6
$ws.Range("H132").Value = 1225264
$ws.Range("L132").Value = 24256.2861
$ws.Range("N132").Value = -29316.2861
$ws.Range("M132").Value = -5566650.800000001
$ws.Range("K132").Value = 5569180.800000001
$ws.Range("J132").Value = 8085.4287
$ws.Range("M136").Value = -6667.600199999999
$ws.Range("I136").Value = 3072.5334
$ws.Range("H136").Value = 7621.121
$ws.Range("K136").Value = 9217.600199999999
$ws.Range("J140").Value = 69796.60000000001
$ws.Range("H140").Value = 69796.60000000001
$ws.Range("L140").Value = 69796.60000000001
$ws.Range("N140").Value = -80156.60000000001
$ws.Range("N141").Value = -35074.5
$ws.Range("H141").Value = 24714.5
$ws.Range("J141").Value = 24714.5
$ws.Range("L141").Value = 24714.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L11").Value = 4366
$ws.Range("I11").Value = 1549.4
$ws.Range("M11").Value = -1409.4
$ws.Range("H11").Value = 2801.2222
$ws.Range("N11").Value = -4646
$ws.Range("K11").Value = 1549.4
$ws.Range("J11").Value = 4366
$ws.Range("H86").Value = 6453573
$ws.Range("K86").Value = 12510514
$ws.Range("I86").Value = 12510514
$ws.Range("M86").Value = -12509391
$ws.Range("H89").Value = 6453573
$ws.Range("M89").Value = -62546954
$ws.Range("I89").Value = 12510514
$ws.Range("K89").Value = 62552570
$ws.Range("H94").Value = 43481004
$ws.Range("I94").Value = 83334500
$ws.Range("K94").Value = 83334500
$ws.Range("M94").Value = -83334049
$ws.Range("N105").Value = -8706.222
$ws.Range("M105").Value = -500.375
$ws.Range("K105").Value = 2247.375
$ws.Range("H105").Value = 3817
$ws.Range("L105").Value = 5212.222
$ws.Range("J105").Value = 5212.222
$ws.Range("I105").Value = 2247.375
$ws.Range("I134").Value = 1942.1428
$ws.Range("K134").Value = 5826.428400000001
$ws.Range("M134").Value = -3291.428400000001
$ws.Range("H134").Value = 6294.1904

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M6").Value = -7744.143
$ws.Range("I6").Value = 7857.143
$ws.Range("H6").Value = 7500
$ws.Range("K6").Value = 7857.143
$ws.Range("I19").Value = 3016.6
$ws.Range("M19").Value = -2846.6
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("H19").Value = 2847.1667
$ws.Range("N19").Value = -2340
$ws.Range("K19").Value = 3016.6
$ws.Range("H22").Value = 320.75
$ws.Range("I22").Value = 295.36365
$ws.Range("M22").Value = 54.63634999999999
$ws.Range("K22").Value = 295.36365
$ws.Range("J24").Value = 2000
$ws.Range("N24").Value = -2340
$ws.Range("H24").Value = 2847.1667
$ws.Range("K24").Value = 3016.6
$ws.Range("M24").Value = -2846.6
$ws.Range("I24").Value = 3016.6
$ws.Range("L24").Value = 2000
$ws.Range("H31").Value = 9143.958000000001
$ws.Range("L31").Value = 11397
$ws.Range("N31").Value = -11987
$ws.Range("J31").Value = 11397
$ws.Range("I32").Value = 2750
$ws.Range("M32").Value = -2434
$ws.Range("H32").Value = 2750
$ws.Range("K32").Value = 2750
$ws.Range("L34").Value = 11397
$ws.Range("N34").Value = -11801
$ws.Range("H34").Value = 9143.958000000001
$ws.Range("J34").Value = 11397
$ws.Range("M62").Value = -17857318
$ws.Range("I62").Value = 17857942
$ws.Range("K62").Value = 17857942
$ws.Range("H62").Value = 15625836
$ws.Range("I65").Value = 17857942
$ws.Range("H65").Value = 15625836
$ws.Range("K65").Value = 89289710
$ws.Range("M65").Value = -89286590
$ws.Range("I132").Value = 7197.4
$ws.Range("H132").Value = 11089.272
$ws.Range("M132").Value = -19062.2
$ws.Range("K132").Value = 21592.2
$ws.Range("I134").Value = 2092.4
$ws.Range("K134").Value = 6277.200000000001
$ws.Range("M134").Value = -3742.200000000001
$ws.Range("H134").Value = 5600.607

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K4").Value = 42472119
$ws.Range("H4").Value = 33234624
$ws.Range("M4").Value = -42472007
$ws.Range("I4").Value = 14157373
$ws.Range("N12").Value = -15000715
$ws.Range("L12").Value = 15000369
$ws.Range("J12").Value = 5000123
$ws.Range("H12").Value = 2272998
$ws.Range("H51").Value = 821.4167
$ws.Range("I68").Value = 1605
$ws.Range("H68").Value = 2370.4348
$ws.Range("M68").Value = -4004
$ws.Range("K68").Value = 4815
$ws.Range("I71").Value = 1605
$ws.Range("M71").Value = -10389
$ws.Range("K71").Value = 14445
$ws.Range("H71").Value = 2370.4348
$ws.Range("M76").Value = -2999999617
$ws.Range("I76").Value = 1000000000
$ws.Range("K76").Value = 3000000000
$ws.Range("H76").Value = 1000000000
$ws.Range("I79").Value = 1000000000
$ws.Range("H79").Value = 1000000000
$ws.Range("M79").Value = -2999998674
$ws.Range("K79").Value = 3000000000
$ws.Range("N121").Value = -11401
$ws.Range("L121").Value = 8781
$ws.Range("H121").Value = 1895.9286
$ws.Range("J121").Value = 2927
$ws.Range("K124").Value = 1500
$ws.Range("M124").Value = 3410
$ws.Range("H124").Value = 500
$ws.Range("I124").Value = 500
$ws.Range("H129").Value = 11178418
$ws.Range("I129").Value = 509.0909
$ws.Range("M129").Value = 3472.7273
$ws.Range("K129").Value = 1527.2727
$ws.Range("I130").Value = 997.5
$ws.Range("H130").Value = 997.5
$ws.Range("M130").Value = 2027.5
$ws.Range("K130").Value = 2992.5
$ws.Range("I132").Value = 4043.4
$ws.Range("H132").Value = 6090.0557
$ws.Range("M132").Value = -33860.6
$ws.Range("K132").Value = 36390.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N43").Value = -12302
$ws.Range("H43").Value = 1655.2222
$ws.Range("L43").Value = 12000
$ws.Range("J43").Value = 12000
$ws.Range("H102").Value = 5723.5293
$ws.Range("I102").Value = 5386.6665
$ws.Range("K102").Value = 5386.6665
$ws.Range("M102").Value = -3764.6665
$ws.Range("I107").Value = 1333523.1
$ws.Range("M107").Value = -1331603.1
$ws.Range("L107").Value = 1163.5714
$ws.Range("N107").Value = -5003.5714
$ws.Range("K107").Value = 1333523.1
$ws.Range("J107").Value = 1163.5714
$ws.Range("H107").Value = 616098.75
$ws.Range("K126").Value = 7503.75
$ws.Range("I126").Value = 2501.25
$ws.Range("H126").Value = 5135.857
$ws.Range("M126").Value = -5033.75
$ws.Range("I132").Value = 2395.4211
$ws.Range("H132").Value = 3941.516
$ws.Range("L132").Value = 19168.5
$ws.Range("N132").Value = -24228.5
$ws.Range("M132").Value = -4656.263300000001
$ws.Range("K132").Value = 7186.263300000001
$ws.Range("J132").Value = 6389.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N7").Value = -8849
$ws.Range("L7").Value = 8625
$ws.Range("H7").Value = 8388.888999999999
$ws.Range("J7").Value = 8625
$ws.Range("M16").Value = -3059.8125
$ws.Range("L16").Value = 3999
$ws.Range("H16").Value = 3351.2632
$ws.Range("J16").Value = 3999
$ws.Range("N16").Value = -4339
$ws.Range("I16").Value = 3229.8125
$ws.Range("K16").Value = 3229.8125
$ws.Range("I40").Value = 4968.2
$ws.Range("K40").Value = 4968.2
$ws.Range("H40").Value = 4968.2
$ws.Range("M40").Value = -4832.2
$ws.Range("K46").Value = 1298.6
$ws.Range("J46").Value = 3928.8572
$ws.Range("L46").Value = 3928.8572
$ws.Range("H46").Value = 2832.9167
$ws.Range("M46").Value = -1110.6
$ws.Range("N46").Value = -4304.8572
$ws.Range("I46").Value = 1298.6
$ws.Range("J61").Value = 6749.875
$ws.Range("I61").Value = 20002460
$ws.Range("L61").Value = 6749.875
$ws.Range("K61").Value = 20002460
$ws.Range("N61").Value = -7153.875
$ws.Range("M61").Value = -20002258
$ws.Range("H61").Value = 7697407.5
$ws.Range("J81").Value = 43000
$ws.Range("H81").Value = 43000
$ws.Range("N81").Value = -44996
$ws.Range("L81").Value = 43000
$ws.Range("H82").Value = 1085324.8
$ws.Range("M82").Value = -1409411.2
$ws.Range("I82").Value = 1409772.2
$ws.Range("K82").Value = 1409772.2
$ws.Range("H84").Value = 43000
$ws.Range("J84").Value = 43000
$ws.Range("N84").Value = -138984
$ws.Range("L84").Value = 129000
$ws.Range("H85").Value = 1085324.8
$ws.Range("I85").Value = 1409772.2
$ws.Range("M85").Value = -1408524.2
$ws.Range("K85").Value = 1409772.2
$ws.Range("H113").Value = 7697407.5
$ws.Range("J113").Value = 6749.875
$ws.Range("L113").Value = 6749.875
$ws.Range("I113").Value = 20002460
$ws.Range("K113").Value = 20002460
$ws.Range("M113").Value = -20000290
$ws.Range("N113").Value = -11089.875
$ws.Range("L126").Value = 25875
$ws.Range("H126").Value = 8388.888999999999
$ws.Range("J126").Value = 8625
$ws.Range("N126").Value = -30815
$ws.Range("M136").Value = -14547.75
$ws.Range("I136").Value = 5699.25
$ws.Range("H136").Value = 8308.76
$ws.Range("K136").Value = 17097.75
$ws.Range("J136").Value = 12947.889
$ws.Range("L136").Value = 38843.667
$ws.Range("N136").Value = -43943.667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N76").Value = -30630
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("H79").Value = 30000
$ws.Range("N79").Value = -32184
$ws.Range("L79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("H113").Value = 21417.584
$ws.Range("I113").Value = 41829.168
$ws.Range("K113").Value = 125487.504
$ws.Range("M113").Value = -123317.504
$ws.Range("L122").Value = 31257.999
$ws.Range("K122").Value = 336930.39
$ws.Range("H122").Value = 78964.05499999999
$ws.Range("I122").Value = 112310.13
$ws.Range("N122").Value = -36157.999
$ws.Range("J122").Value = 10419.333
$ws.Range("M122").Value = -334480.39
$ws.Range("J123").Value = 51264
$ws.Range("L123").Value = 51264
$ws.Range("H123").Value = 51264
$ws.Range("N123").Value = -61064
$ws.Range("K126").Value = 5727.6
$ws.Range("I126").Value = 1909.2
$ws.Range("H126").Value = 4241.737
$ws.Range("M126").Value = -3257.6
$ws.Range("M136").Value = -1701.642599999999
$ws.Range("I136").Value = 1417.2142
$ws.Range("H136").Value = 27184.906
$ws.Range("K136").Value = 4251.642599999999
$ws.Range("J136").Value = 75284.60000000001
$ws.Range("L136").Value = 225853.8
$ws.Range("N136").Value = -230953.8
